$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The class-schedule row for Friday 13:00-14:00 (ว 42002) is relabelled
# from ม.6/1 to ม.6/2.
$ws.Range("C7").Value = "ม.6/2"

# Selected cell in the sheet view moves from H5 to J13.
$ws.Range("J13").Select()
